$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cellule A1 vide"
$ws.Range("A2").Value = "Cellule A2"

$ws.Range("A2").Select()
